# Added results click tests
#
# 1. Reorder tabs: DictionaryPopup now comes before SearchTerms.
# 2. Definitions becomes the active/selected tab (was SearchTerms before).
# 3. SearchTerms sheet gets new "results click test" data: for each
#    dictionary, a search term / search type / link name used by the new
#    tests.
# 4. TermsPageDrugs sheet's remembered selection moves to A5.

$wb = $excel.ActiveWorkbook

# --- 1. Move DictionaryPopup in front of SearchTerms ---------------------
$dictionaryPopup = $wb.Worksheets.Item("DictionaryPopup")
$searchTermsBeforeMove = $wb.Worksheets.Item("SearchTerms")
$dictionaryPopup.Move($searchTermsBeforeMove)

# Re-resolve by name: after Move(), old references track sheet *position*,
# not identity, so grab a fresh handle to the sheet that is now named
# "SearchTerms".
$searchTerms = $wb.Worksheets.Item("SearchTerms")

# --- 2. Replace SearchTerms contents with the new click-test data --------
$searchTerms.Cells.Item(2, 1).Value = "/publications/dictionaries/cancer-terms"
$searchTerms.Cells.Item(2, 2).Value = "breast"
$searchTerms.Cells.Item(2, 3).Value = "dictionary_terms"
$searchTerms.Cells.Item(2, 4).Value = "TermsDictionarySearch"

$searchTerms.Cells.Item(3, 1).Value = "/espanol/publicaciones/diccionario"
$searchTerms.Cells.Item(3, 2).Value = "tumor"
$searchTerms.Cells.Item(3, 3).Value = "diccionario"
$searchTerms.Cells.Item(3, 4).Value = "TermsDictionarySearch"

$searchTerms.Cells.Item(4, 1).Value = "/publications/dictionaries/cancer-drug"
$searchTerms.Cells.Item(4, 2).Value = "herceptin"
$searchTerms.Cells.Item(4, 3).Value = "dictionary_drugs"
$searchTerms.Cells.Item(4, 4).Value = "DrugDictionarySearch"

$searchTerms.Cells.Item(5, 1).Value = "/publications/dictionaries/genetics-dictionary"
$searchTerms.Cells.Item(5, 2).Value = "allele"
$searchTerms.Cells.Item(5, 3).Value = "dictionary_genetics"
$searchTerms.Cells.Item(5, 4).Value = "GeneticsDictionarySearch"

# Columns got wider once they picked up SearchType/LinkName values.
$searchTerms.Columns.Item(1).ColumnWidth = 43
$searchTerms.Columns.Item(2).ColumnWidth = 11.42578125
$searchTerms.Columns.Item(3).ColumnWidth = 18.5703125
$searchTerms.Columns.Item(4).ColumnWidth = 24

[void]$searchTerms.Range("A6").Select()

# --- 3. TermsPageDrugs: remembered selection moves to A5 ------------------
$termsPageDrugs = $wb.Worksheets.Item("TermsPageDrugs")
[void]$termsPageDrugs.Range("A5").Select()

# --- 4. Definitions becomes the active tab --------------------------------
$wb.Worksheets.Item("Definitions").Activate() | Out-Null
[void]$wb.Worksheets.Item("Definitions").Range("A10").Select()
